# chore: adapt column header formatting to respective input file names
#
# - Renames the "_old" / "_new" header-suffix convention used in the
#   AHB-Diff sheet's header row to the concrete format-version names
#   "_FV2404" (old/left side) and "_FV2410" (new/right side). The lone
#   "diff" header in between is left untouched.
# - Wraps the used range in an Excel Table ("Table1") with a header row
#   + autofilter so the table shows up as a first-class ListObject.
# - Freezes the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
$oldSuffix = "_old"
$newSuffixOld = "_FV2404"
$suffixNew = "_new"
$newSuffixNew = "_FV2410"

$lastCol = 21   # A:U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $headerText = $cell.Value2

    if ($headerText -like "*$oldSuffix") {
        $base = $headerText.Substring(0, $headerText.Length - $oldSuffix.Length)
        $cell.Value = $base + $newSuffixOld
    } elseif ($headerText -like "*$suffixNew") {
        $base = $headerText.Substring(0, $headerText.Length - $suffixNew.Length)
        $cell.Value = $base + $newSuffixNew
    }
}

# --- 2. Turn the used range into an Excel Table ------------------------
$tableRange = $ws.Range("A1:U52")
$table = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
